# Insert a new data row before the existing row 52 (Excel shifts rows 52..130
# down to 53..131, updates the used-range dimension automatically).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(52).Insert()

# Populate the newly-inserted row 52 with the new "Haba" price entry.
$ws.Cells.Item(52, 1).Value = 4
$ws.Cells.Item(52, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(52, 3).Value = "Los Lagos"
$ws.Cells.Item(52, 4).Value = 45079
$ws.Cells.Item(52, 5).Value = 10
$ws.Cells.Item(52, 6).Value = 100112026
$ws.Cells.Item(52, 7).Value = "Haba"
$ws.Cells.Item(52, 8).Value = "Sin especificar"
$ws.Cells.Item(52, 9).Value = "Primera"
$ws.Cells.Item(52, 10).Value = 80
$ws.Cells.Item(52, 11).Value = 21000
$ws.Cells.Item(52, 12).Value = 21000
$ws.Cells.Item(52, 13).Value = 21000
$ws.Cells.Item(52, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(52, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(52, 16).Value = 840
$ws.Cells.Item(52, 17).Value = 25
$ws.Cells.Item(52, 18).Value = "Hortaliza"
